# Apply edit: add "metadata" sheet, refresh F-column (time_taken) timestamps on "data" sheet

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the time_taken values (column F, rows 2..66) on the "data" sheet ---
$newTimestamps = @(
    "2021-10-05 14:22:56.193003",
    "2021-10-05 14:22:56.193012",
    "2021-10-05 14:22:56.193016",
    "2021-10-05 14:22:56.193019",
    "2021-10-05 14:22:56.193022",
    "2021-10-05 14:22:56.193025",
    "2021-10-05 14:22:56.193027",
    "2021-10-05 14:22:56.193030",
    "2021-10-05 14:22:56.193033",
    "2021-10-05 14:22:56.193036",
    "2021-10-05 14:22:56.193038",
    "2021-10-05 14:22:56.193041",
    "2021-10-05 14:22:56.193044",
    "2021-10-05 14:22:56.193047",
    "2021-10-05 14:22:56.193049",
    "2021-10-05 14:22:56.193052",
    "2021-10-05 14:22:56.193055",
    "2021-10-05 14:22:56.193058",
    "2021-10-05 14:22:56.193061",
    "2021-10-05 14:22:56.193063",
    "2021-10-05 14:22:56.193066",
    "2021-10-05 14:22:56.193069",
    "2021-10-05 14:22:56.193072",
    "2021-10-05 14:22:56.193074",
    "2021-10-05 14:22:56.193077",
    "2021-10-05 14:22:56.193080",
    "2021-10-05 14:22:56.193083",
    "2021-10-05 14:22:56.193086",
    "2021-10-05 14:22:56.193088",
    "2021-10-05 14:22:56.193091",
    "2021-10-05 14:22:56.193094",
    "2021-10-05 14:22:56.193097",
    "2021-10-05 14:22:56.193100",
    "2021-10-05 14:22:56.193103",
    "2021-10-05 14:22:56.193106",
    "2021-10-05 14:22:56.193108",
    "2021-10-05 14:22:56.193111",
    "2021-10-05 14:22:56.193114",
    "2021-10-05 14:22:56.193116",
    "2021-10-05 14:22:56.193119",
    "2021-10-05 14:22:56.193122",
    "2021-10-05 14:22:56.193125",
    "2021-10-05 14:22:56.193127",
    "2021-10-05 14:22:56.193130",
    "2021-10-05 14:22:56.193133",
    "2021-10-05 14:22:56.193135",
    "2021-10-05 14:22:56.193138",
    "2021-10-05 14:22:56.193141",
    "2021-10-05 14:22:56.193143",
    "2021-10-05 14:22:56.193146",
    "2021-10-05 14:22:56.193149",
    "2021-10-05 14:22:56.193152",
    "2021-10-05 14:22:56.193155",
    "2021-10-05 14:22:56.193158",
    "2021-10-05 14:22:56.193161",
    "2021-10-05 14:22:56.193163",
    "2021-10-05 14:22:56.193166",
    "2021-10-05 14:22:56.193169",
    "2021-10-05 14:22:56.193172",
    "2021-10-05 14:22:56.193174",
    "2021-10-05 14:22:56.193177",
    "2021-10-05 14:22:56.193180",
    "2021-10-05 14:22:56.193183",
    "2021-10-05 14:22:56.193186",
    "2021-10-05 14:22:56.193190"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- 2. Add the "metadata" sheet after "data" ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Reuse the header style/formatting from the "data" sheet by copying cells, then overwrite text
$dataSheet.Range("B1:F1").Copy($metaSheet.Range("B1:F1"))
$dataSheet.Range("F1").Copy($metaSheet.Range("G1"))
$dataSheet.Range("A2").Copy($metaSheet.Range("A2"))

$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

$metaSheet.Cells.Item(2, 2).Value = "Thoracic aortic aneurysm or dissection"
$metaSheet.Cells.Item(2, 3).Value = 1
$metaSheet.Cells.Item(2, 4).Value = "'1.123"
$metaSheet.Cells.Item(2, 5).Value = "2021-09-14T10:18:09.378495Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:22:56.189352"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/1/?format=json"

$dataSheet.Activate()
$dataSheet.Range("A1").Select()
